$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text (avoid Excel auto-converting
# numeric-looking strings like "241.93" or "1.002" into real numbers),
# then restore the default "Normal" style so no stray style/format is left
# behind on the cell (matches original formatting exactly).
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "30.795.18"
Set-TextValue $ws.Range("E2") "  -1.40%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.945.22"
Set-TextValue $ws.Range("E3") "  -1.18%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.32%  "

# Row 5
Set-TextValue $ws.Range("D5") "241.93"
Set-TextValue $ws.Range("E5") "  -2.54%  "

# Row 6
Set-TextValue $ws.Range("E6") "  +0.19%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4860"
Set-TextValue $ws.Range("E7") "  -0.65%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.2920"
Set-TextValue $ws.Range("E8") "  -1.97%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.06845"
Set-TextValue $ws.Range("E9") "  -0.23%  "

# Row 10
Set-TextValue $ws.Range("D10") "19.42"
Set-TextValue $ws.Range("E10") "  +0.73%  "

# Row 11
Set-TextValue $ws.Range("D11") "105.17"
Set-TextValue $ws.Range("E11") "  -1.80%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.956.50"
Set-TextValue $ws.Range("E12") "  +0.36%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.07744"
Set-TextValue $ws.Range("E13") "  -0.47%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.303"
Set-TextValue $ws.Range("E14") "  -2.91%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.6943"
Set-TextValue $ws.Range("E15") "  -3.49%  "

# Row 16
Set-TextValue $ws.Range("D16") "273.90"
Set-TextValue $ws.Range("E16") "  -4.27%  "

# Row 17
Set-TextValue $ws.Range("D17") "30.822.37"
Set-TextValue $ws.Range("E17") "  -0.92%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.000007690"
Set-TextValue $ws.Range("E18") "  -1.08%  "

# Row 19
Set-TextValue $ws.Range("B19") "Avalanche"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "13.10"
Set-TextValue $ws.Range("E19") "  -1.84%  "

# Row 20
Set-TextValue $ws.Range("B20") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D20") "2.201.75"
Set-TextValue $ws.Range("E20") "  -0.51%  "

# Row 21
Set-TextValue $ws.Range("E21") "  +0.09%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.440"
Set-TextValue $ws.Range("E22") "  -3.86%  "

# Row 23
Set-TextValue $ws.Range("D23") "1.002"
Set-TextValue $ws.Range("E23") "  +0.20%  "

# Row 24
Set-TextValue $ws.Range("D24") "6.450"
Set-TextValue $ws.Range("E24") "  -3.26%  "

# Row 25
Set-TextValue $ws.Range("D25") "9.680"
Set-TextValue $ws.Range("E25") "  -4.06%  "

# Row 26
Set-TextValue $ws.Range("D26") "167.26"
Set-TextValue $ws.Range("E26") "  -1.25%  "

# Row 27
Set-TextValue $ws.Range("D27") "19.49"
Set-TextValue $ws.Range("E27") "  -2.89%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.156"
Set-TextValue $ws.Range("E28") "  -2.66%  "

# Row 29
Set-TextValue $ws.Range("E29") "  -3.25%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.392"
Set-TextValue $ws.Range("E30") "  -3.83%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.553"
Set-TextValue $ws.Range("E31") "  -2.92%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.513"
Set-TextValue $ws.Range("E32") "  -7.40%  "

# Row 33
Set-TextValue $ws.Range("D33") "4.357"
Set-TextValue $ws.Range("E33") "  -4.05%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.04831"
Set-TextValue $ws.Range("E34") "  -5.20%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.7428"
Set-TextValue $ws.Range("E35") "  -3.90%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.152"
Set-TextValue $ws.Range("E36") "  -2.12%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.730"
Set-TextValue $ws.Range("E37") "  -0.14%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.01981"
Set-TextValue $ws.Range("E38") "  -3.84%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.667"
Set-TextValue $ws.Range("E39") "  -1.66%  "

# Row 40
Set-TextValue $ws.Range("D40") "6.447"
Set-TextValue $ws.Range("E40") "  -0.01%  "

# Row 41
Set-TextValue $ws.Range("D41") "76.71"
Set-TextValue $ws.Range("E41") "  +3.84%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.081"
Set-TextValue $ws.Range("E42") "  -2.54%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.8940"
Set-TextValue $ws.Range("E43") "  +0.75%  "

# Row 44
Set-TextValue $ws.Range("D44") "107.92"
Set-TextValue $ws.Range("E44") "  -2.02%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.4391"
Set-TextValue $ws.Range("E45") "  -2.36%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.9990"
Set-TextValue $ws.Range("E46") "  -0.15%  "

# Row 47
Set-TextValue $ws.Range("D47") "7.701"
Set-TextValue $ws.Range("E47") "  +2.29%  "

# Row 48
Set-TextValue $ws.Range("D48") "996.61"
Set-TextValue $ws.Range("E48") "  -0.07%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.1236"
Set-TextValue $ws.Range("E49") "  -2.96%  "

# Row 50
Set-TextValue $ws.Range("B50") "Elrond"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D50") "35.54"
Set-TextValue $ws.Range("E50") "  -1.64%  "

# Row 51
Set-TextValue $ws.Range("B51") "EnergySwap"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "9.092"
Set-TextValue $ws.Range("E51") "  -3.99%  "
